$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GET_activity_v1_users")
$ws2 = $wb.Worksheets.Item("GET_last_login")

# Add the new test-case row (row 3) on the GET_last_login sheet
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Y"
$ws2.Range("C3").Value = "Get Last login with invalid userId"
$ws2.Range("E3").Value = "GET"
$ws2.Range("F3").Value = "/activity/v1/users/a416e744-c66f-48b1-af78-055aa30aa982/last-logins"
$ws2.Range("G3").Value = "400"

# Update selections: first sheet is no longer the active/selected tab,
# second sheet becomes active with its own new selection.
$ws1.Range("F8").Select()
$ws2.Activate()
$ws2.Range("G3").Select()
